$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "dwadwaBestand 1 column 4"
$ws.Range("C6").Value = "dadaBestand 1 column 6"
$ws.Range("C7").Value = "adaBestand 1 column 2"
$ws.Range("C9").Value = "dadaBestand 1 column 8"
$ws.Range("C10").Value = "dadaBestand 1 column 11"
$ws.Range("C11").Value = "dadaBestand 1 column 10"
$ws.Range("C12").Value = "adaBestand 1 column 9"
